$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new "MISC Change" column before the existing "FACE Change" column (B) ---
$ws.Columns(2).Insert()
$ws.Columns(2).ColumnWidth = $ws.Columns(1).ColumnWidth
$ws.Cells.Item(1, 2).Value = "MISC Change"

# --- 2. Fill in the new release row (row 5) that was previously blank ---
$ws.Range("A5").Value = "CHECON.MECH.v.1.3"
$ws.Range("B5").Value = "Added BOM"
$ws.Range("C5").Value = "None"
$ws.Range("D5").Value = "None"
$ws.Range("E5").Value = "None"
$ws.Range("F5").Value = "None"
$ws.Range("G5").Value = "None"
$ws.Range("H5").Value = "None"
$ws.Range("I5").Value = "None"
$ws.Range("J5").Value = "Ethan Dale"
$ws.Range("K5").Value = "Subassembly requires its own BOM (value add)"

# Give row 5's text cells the same look (font/border/alignment) as the row above it.
$ws.Range("A4:K4").Copy()
$ws.Range("A5:K5").PasteSpecial(-4122)
$ws.Range("A5").Value = "CHECON.MECH.v.1.3"
$ws.Range("B5").Value = "Added BOM"
$ws.Range("C5").Value = "None"
$ws.Range("D5").Value = "None"
$ws.Range("E5").Value = "None"
$ws.Range("F5").Value = "None"
$ws.Range("G5").Value = "None"
$ws.Range("H5").Value = "None"
$ws.Range("I5").Value = "None"
$ws.Range("J5").Value = "Ethan Dale"
$ws.Range("K5").Value = "Subassembly requires its own BOM (value add)"

# Date of the new release, formatted m/d/yyyy (matches the new date style used for L4/L5)
$ws.Range("K2").Copy()
$ws.Range("L5").PasteSpecial(-4122)
$ws.Range("L5").Value = 42129
$ws.Range("L5").NumberFormat = "m/d/yyyy"

# --- 3. Record a release date for the v1.2 row (L4), which had none before ---
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("L4").Value = 42116
$ws.Range("L4").NumberFormat = "m/d/yyyy"
